$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.899.24'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '1.668.26'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = "'215.76"
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = "'0.535"
$ws.Range('E6').Value = '  +5.42%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +1.21%  '
$ws.Range('E9').Value = '  +0.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = "'20.26"
$ws.Range('E10').Value = '  +2.93%  '
$ws.Range('E11').Value = '  +3.46%  '
$ws.Range('D12').Value = '1.903.68'
$ws.Range('E12').Value = '  +1.15%  '
$ws.Range('D13').Value = '1.647.70'
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('E14').Value = '  +0.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = "'0.527"
$ws.Range('E15').Value = '  +1.84%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = "'66.13"
$ws.Range('E16').Value = '  +1.52%  '
$ws.Range('D17').Value = '26.941.55'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = "'234.25"
$ws.Range('E18').Value = '  -1.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = "'7.94"
$ws.Range('E19').Value = '  +1.55%  '
$ws.Range('E20').Value = '  +0.53%  '
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = "'4.39"
$ws.Range('E22').Value = '  -0.32%  '
$ws.Range('B23').Value = 'Avalanche'
$ws.Range('C23').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = "'9.14"
$ws.Range('E23').Value = '  -0.97%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = "'2.19"
$ws.Range('E24').Value = '  -2.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = "'146.20"
$ws.Range('E25').Value = '  +0.43%  '
$ws.Range('E26').Value = '  +0.72%  '
$ws.Range('E27').Value = '  +1.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = "'15.92"
$ws.Range('E28').Value = '  +0.74%  '
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = "'0.0497"
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('E31').Value = '  +0.03%  '
$ws.Range('E32').Value = '  +2.09%  '
$ws.Range('D33').Value = '1.453.76'
$ws.Range('E33').Value = '  -3.76%  '
$ws.Range('E34').Value = '  +2.34%  '
$ws.Range('E35').Value = '  +4.09%  '
$ws.Range('E36').Value = '  -0.39%  '
$ws.Range('E37').Value = '  +1.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = "'0.903"
$ws.Range('E38').Value = '  +2.13%  '
$ws.Range('E39').Value = '  +0.68%  '
$ws.Range('E40').Value = '  -3.40%  '
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('E42').Value = '  +1.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = "'66.15"
$ws.Range('E43').Value = '  +0.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = "'0.973"
$ws.Range('E44').Value = '  +6.18%  '
$ws.Range('D45').Value = '1.812.33'
$ws.Range('E45').Value = '  +1.24%  '
$ws.Range('E46').Value = '  +1.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = "'90.65"
$ws.Range('E47').Value = '  +1.44%  '
$ws.Range('E48').Value = '  +1.56%  '
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('E50').Value = '  +4.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = "'0.0506"
$ws.Range('E51').Value = '  -0.17%  '
